$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.711.26'
$ws.Range("E2").Value = '  +0.01%  '

# Row 3
$ws.Range("D3").Value = '3.511.88'
$ws.Range("E3").Value = '  +0.22%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.83'
$ws.Range("E5").Value = '  -1.06%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '195.49'
$ws.Range("E6").Value = '  +2.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  -0.96%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.202'
$ws.Range("E9").Value = '  -5.31%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.643'
$ws.Range("E10").Value = '  -2.83%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.13'
$ws.Range("E11").Value = '  -0.40%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000301'
$ws.Range("E12").Value = '  -1.89%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.46'
$ws.Range("E13").Value = '  -1.29%  '

# Row 14
$ws.Range("D14").Value = '4.067.16'
$ws.Range("E14").Value = '  +0.10%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '600.59'
$ws.Range("E15").Value = '  -2.34%  '

# Row 16
$ws.Range("D16").Value = '69.806.33'
$ws.Range("E16").Value = '  -0.02%  '

# Row 17
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.70'
$ws.Range("E17").Value = '  +0.67%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.94'
$ws.Range("E18").Value = '  +0.06%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.122'
$ws.Range("E19").Value = '  +1.99%  '

# Row 20
$ws.Range("D20").Value = '3.510.43'
$ws.Range("E20").Value = '  +0.20%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.981'
$ws.Range("E21").Value = '  -0.65%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.07'
$ws.Range("E22").Value = '  +4.34%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.20'
$ws.Range("E23").Value = '  +2.66%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.65'
$ws.Range("E24").Value = '  -3.16%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.60'
$ws.Range("E25").Value = '  -1.00%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.06'
$ws.Range("E26").Value = '  -0.33%  '

# Row 27
$ws.Range("E27").Value = '  -1.59%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.49'
$ws.Range("E28").Value = '  -2.06%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.05'
$ws.Range("E29").Value = '  -4.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.37'
$ws.Range("E30").Value = '  +8.02%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.98'
$ws.Range("E31").Value = '  +0.40%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.28'
$ws.Range("E32").Value = '  -1.06%  '

# Row 33
$ws.Range("E33").Value = '  -1.32%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.10'
$ws.Range("E34").Value = '  -1.04%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.16'
$ws.Range("E35").Value = '  +2.38%  '

# Row 36
$ws.Range("D36").Value = '3.745.60'
$ws.Range("E36").Value = '  +1.91%  '

# Row 37
$ws.Range("D37").Value = '0.0₃0818'
$ws.Range("E37").Value = '  +5.82%  '

# Row 38
$ws.Range("E38").Value = '  +0.14%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.59'
$ws.Range("E39").Value = '  +0.09%  '

# Row 40
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '501.80'
$ws.Range("E40").Value = '  -3.96%  '

# Row 41
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.388'
$ws.Range("E41").Value = '  -0.61%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.98'
$ws.Range("E42").Value = '  -1.98%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.133'
$ws.Range("E43").Value = '  -3.56%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0448'
$ws.Range("E44").Value = '  -3.39%  '

# Row 45
$ws.Range("E45").Value = '  -1.97%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("E46").Value = '  -2.75%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.26'
$ws.Range("E47").Value = '  -1.91%  '

# Row 48
$ws.Range("E48").Value = '  +0.06%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.43'
$ws.Range("E49").Value = '  -3.31%  '

# Row 50
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000246'
$ws.Range("E50").Value = '  +3.21%  '

# Row 51
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.32'
$ws.Range("E51").Value = '  +3.41%  '
